# ============================================================================
# Edit script: add a new "Expected Label"/"Expected Result" column (U) to the
# FastQuotesConfig sheet, re-color the header/data rows, and tweak a couple
# of alignment / selection details.
# ============================================================================

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ----------------------------------------------------------------------
# 1. New column U: width + content mirroring column T ("Expected" block)
# ----------------------------------------------------------------------
$ws.Columns("U").ColumnWidth = 21.1

$ws.Range("U1").Value = $null
$ws.Range("U2").Value = $null
$ws.Range("U3").Value = "Expected Label"
$ws.Range("U4").Value = $ws.Range("T4").Value

# Merge the new "Expected" header/sub-header cells with column T, mirroring
# the existing T1/T2 merges.
$ws.Range("T1:U1").Merge() | Out-Null
$ws.Range("T2:U2").Merge() | Out-Null

# ----------------------------------------------------------------------
# 2. Formatting: copy the look & feel of the neighbouring cells onto the
#    new column so the new cells blend in with the existing rows.
# ----------------------------------------------------------------------
$ws.Range("T1").Copy() | Out-Null
$ws.Range("U1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("T2").Copy() | Out-Null
$ws.Range("U2").PasteSpecial(-4122) | Out-Null

$ws.Range("T3").Copy() | Out-Null
$ws.Range("U3").PasteSpecial(-4122) | Out-Null

$ws.Range("T4").Copy() | Out-Null
$ws.Range("U4").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# Center-align the merged "Expected" title in T1:U1 (was left aligned).
$ws.Range("T1:U1").HorizontalAlignment = -4108   # xlCenter

# Give T2:U2 the same border as the rest of row 2 (drop the special
# "no right border" style that only T2 used to have).
$ws.Range("R2").Copy() | Out-Null
$ws.Range("T2").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("T2:U2").HorizontalAlignment = -4108   # xlCenter
$ws.Range("T2:U2").Interior.ThemeColor = 6
$ws.Range("T2:U2").Font.Color = 0

# ----------------------------------------------------------------------
# 3. Highlight the column-header row (row 3) and the data row (row 4)
#    with new fill colors (Accent3 Lighter 60% / Accent4 Lighter 80%).
# ----------------------------------------------------------------------
$headerRow = $ws.Range("A3:U3")
$headerRow.Interior.ThemeColor = 7
$headerRow.Interior.TintAndShade = 0.6

$dataRow = $ws.Range("B4:R4")
$dataRow.Interior.ThemeColor = 8
$dataRow.Interior.TintAndShade = 0.8
$dataRow.Font.ThemeColor = 1

$ws.Range("U4").Interior.Color = $ws.Range("T4").Interior.Color

# ----------------------------------------------------------------------
# 4. Selection / view tweaks: select A1:S1 (no explicit active cell) and
#    drop the old "scrolled to column M" view state.
# ----------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("A1:S1").Select() | Out-Null
